$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("measRates")
$ws.Activate()
$ws.Range("L31").Select()
$sel = $excel.Selection
Write-Output ($sel.Areas | Get-Member | Out-String)
